$d = $word.ActiveDocument

function Find-ParagraphIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return 0
}

# ------------------------------------------------------------------
# 1) Insert a new list item "Se rends sur la page de modification de
#    groupe" right before the existing "Affiche la page de
#    modification d'un groupe de discussion" item, carrying the
#    _GoBack bookmark that used to sit further down in the document.
# ------------------------------------------------------------------
$targetIdx = Find-ParagraphIndex("Affiche la page de modification d")
if ($targetIdx -eq 0) {
    throw "Could not find target paragraph for insertion"
}

# Insert an empty paragraph before the target, inheriting its list
# formatting, then fill it with the new sentence.
$d.Paragraphs($targetIdx).Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($targetIdx)
$newStart = $newPara.Range.Start
$newPara.Range.Text = "Se rends sur la page de modification de groupe"

# Bookmark the freshly typed text (stop one char short of the
# paragraph's own range end so the mark stays inside this paragraph
# rather than bleeding into the following one).
$newPara2 = $d.Paragraphs($targetIdx)
$bmRange = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2) Remove the old _GoBack bookmark that used to sit around the
#    "-" / " « ERR_USER_ALREADY_INVITED »" boundary.
# ------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
if ($oldBm.Start -ne $newStart) {
    $oldBm.Delete()
}

# ------------------------------------------------------------------
# 3) Merge the "FORMER" and "_" runs into a single "FORMER_" run
#    without disturbing the neighbouring "« ERR_" / "INVALID_
#    INVITATION »" runs (which share the same, empty, formatting and
#    would otherwise get swept into the rebuild too).
# ------------------------------------------------------------------
$former = $d.Content.Duplicate
$former.Find.Execute("FORMER_")
if (-not $former.Find.Found) {
    throw "Could not find FORMER_ run"
}
$fStart = $former.Start
$fEnd = $former.End

$guard = $d.Range($fStart, $fEnd)
$guard.Font.Bold = $true
$rebuild = $d.Range($fStart, $fEnd)
$rebuild.Text = "FORMER_X"
$extra = $d.Range($fEnd, $fEnd + 1)
$extra.Delete()
$merged = $d.Range($fStart, $fEnd)
$merged.Font.Bold = $false
